$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(1).ColumnWidth = 17.77734375
$ws.Columns.Item(2).ColumnWidth = 17.77734375
$ws.Columns.Item(3).ColumnWidth = 20.77734375

# Header row
$ws.Range("A1").Value = "Client's First Name"
$ws.Range("B1").Value = "Client's Last Name"
$ws.Range("C1").Value = "Clinet's Email Address"
$ws.Range("A1:C1").Style = "Input"

# Data row
$ws.Range("A2").Value = "Thomas"
$ws.Range("B2").Value = "Evans"
$ws.Range("C2").Value = "tsevans@vt.edu"

# Hyperlink (mail) on the email cell
$null = $ws.Hyperlinks.Add($ws.Range("C2"), "mailto:tsevans@vt.edu")

# Select the last-edited cell, matching author's saved selection
$null = $ws.Range("C2").Select()
